$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 392 and 393, shifting the existing rows (392..443) down to (394..445).
$ws.Rows.Item(392).Insert()
$ws.Rows.Item(393).Insert()

# Populate the newly inserted row 392 with its data.
$ws.Cells.Item(392, 1).Value = 9
$ws.Cells.Item(392, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(392, 3).Value = "Metropolitana"
$ws.Cells.Item(392, 4).Value = 44769
$ws.Cells.Item(392, 5).Value = 13
$ws.Cells.Item(392, 6).Value = 100112013
$ws.Cells.Item(392, 7).Value = "Alcachofa"
$ws.Cells.Item(392, 8).Value = "Española"
$ws.Cells.Item(392, 9).Value = "Extra"
$ws.Cells.Item(392, 10).Value = 34
$ws.Cells.Item(392, 11).Value = 20000
$ws.Cells.Item(392, 12).Value = 20000
$ws.Cells.Item(392, 13).Value = 20000
$ws.Cells.Item(392, 14).Value = "$/caja 25 unidades"
$ws.Cells.Item(392, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(392, 16).Value = 20000
$ws.Cells.Item(392, 17).Value = 1
$ws.Cells.Item(392, 18).Value = "Hortaliza"

# Populate the newly inserted row 393 with its data.
$ws.Cells.Item(393, 1).Value = 9
$ws.Cells.Item(393, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(393, 3).Value = "Metropolitana"
$ws.Cells.Item(393, 4).Value = 44769
$ws.Cells.Item(393, 5).Value = 13
$ws.Cells.Item(393, 6).Value = 100112013
$ws.Cells.Item(393, 7).Value = "Alcachofa"
$ws.Cells.Item(393, 8).Value = "Española"
$ws.Cells.Item(393, 9).Value = "Primera"
$ws.Cells.Item(393, 10).Value = 52
$ws.Cells.Item(393, 11).Value = 17000
$ws.Cells.Item(393, 12).Value = 17000
$ws.Cells.Item(393, 13).Value = 17000
$ws.Cells.Item(393, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(393, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(393, 16).Value = 567
$ws.Cells.Item(393, 17).Value = 30
$ws.Cells.Item(393, 18).Value = "Hortaliza"

# Ensure the date cells use the same date-formatted number format as the rest of column D.
$ws.Cells.Item(392, 4).NumberFormat = $ws.Cells.Item(391, 4).NumberFormat
$ws.Cells.Item(393, 4).NumberFormat = $ws.Cells.Item(391, 4).NumberFormat
